# Arreglo de importar excel y botones en crear producto
# - Rename the "unidadmedida" header to "clase" (column L header)
# - Widen column C (nombre) so longer product names are readable
# - Shrink column L now that it holds the shorter "clase" label instead of "unidadmedida"
# - Move the active selection to I5 (where the user was last working)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row relabeling: "unidadmedida" -> "clase"
$ws.Range("L1").Value = "clase"

# Column width adjustments
$ws.Columns.Item(3).ColumnWidth = 33      # C: nombre - widen for long product names
$ws.Columns.Item(12).ColumnWidth = 6.5    # L: clase - narrower than old unidadmedida column

# Restore selection to I5
$null = $ws.Range("I5").Select()
